$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Documentation" header to "Documentation/etc"
$ws.Range("G1").Value = "Documentation/etc"

# Widen column G to fit new header text
$ws.Range("G1").EntireColumn.ColumnWidth = 22.666666666666668

# Update formulas that now include an extra entry
$ws.Range("G3").Formula = "=(1/60)*(180+20+10+6)"
$ws.Range("E27").Formula = "=(1/60)*(21+26+20+20)"

# Move the active selection to H12
$ws.Range("H12").Select()
